# Update column G (header "K") values for rows 2-9 per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 3
$ws.Range("G9").Value = 2
